$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Stuff"
$ws.Range("A2").Value = "May have done a bad design change, due to wanting to change values within movement modifiers during runtime such as gravityMagnitude, groundPullMagnitude, forceDrag, movementSpeed, jumpForce, etc…while keeping the movement modifier objects as POCO's without the SerializeField attributes, and the class itself not being serializable would mean that the constructor wouldn't get called automatically each time in the editor, to achieve this the values were delegated over to the corresponding MonoBehaviours and are passed in for the .Tick methods`nBut this means that data builders can't work for them, and that they'll have to be passed as parameters each time. Perhaps this is okay with the current way it's set up but it doesn't sound exactly great`nThe methods should have parameters that affect them, but they should function as lone objects as well, and this might require modification which in turn breaks the open closed principle as well`nPerhaps change this but keep the classes as POCO? But then the values would lose the ability to be changed during runtime in the inspector...?"

$ws.Range("A1:A2").WrapText = $true

$ws.Columns.Item(1).ColumnWidth = 79.5

$ws.Rows.Item(2).RowHeight = 225

$ws.Range("A2").Select()
